$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "nemad" (ticker symbol) column L previously duplicated the
# company name ("نفت سپاهان"); update it to the actual symbol "شسپا"
# for every data row (rows 2-45). This introduces a new shared string.
for ($r = 2; $r -le 45; $r++) {
    $ws.Range("L$r").Value = "شسپا"
}

# Column L was resized (best-fit) to width 10.
$ws.Columns("L").ColumnWidth = 9.166666666666666

# The active selection when the file was last saved was L8.
$ws.Range("L8").Select() | Out-Null
